# Enhance config file handling
# Append a new row (row 61) of config/ID data to each of the 4 sheets,
# mirroring the structure of the existing rows (time, lengths, IDs, checksums).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = 1
        A = 45847.46273148148
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x60"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 352
        I = 7
    },
    @{
        Sheet = 2
        A = 45847.46273148148
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x5C"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 348
        I = 25
    },
    @{
        Sheet = 3
        A = 45847.46273148148
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x68"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 104
        I = 15
    },
    @{
        Sheet = 4
        A = 45847.46273148148
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7D"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 125
        I = 9
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = 61

    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
